# Selenium Season 3 Nov 1st
# Target sheet is "com.adactin" (the active/selected sheet in the workbook).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: make room for a new "Expected Title" column at B4 ---
# Shift C4:G4 one column to the right (into D4:H4). Must go right-to-left,
# cell by cell, so a source cell is never clobbered before it's read.
# Range.Copy(dest) carries both the value and the cell style/format.
$ws.Range("G4").Copy($ws.Range("H4"))
$ws.Range("F4").Copy($ws.Range("G4"))
$ws.Range("E4").Copy($ws.Range("F4"))
$ws.Range("D4").Copy($ws.Range("E4"))
$ws.Range("C4").Copy($ws.Range("D4"))
$ws.Range("B4").Copy($ws.Range("C4"))
$ws.Range("B4").Value2 = "Expected Title"

# --- Row 5: same column shift, then fill the new B5/C5 with the expected
# title and the new username for the registration test case ---
$ws.Range("G5").Copy($ws.Range("H5"))
$ws.Range("F5").Copy($ws.Range("G5"))
$ws.Range("E5").Copy($ws.Range("F5"))
$ws.Range("D5").Copy($ws.Range("E5"))
$ws.Range("C5").Copy($ws.Range("D5"))
$ws.Range("B5").Copy($ws.Range("C5"))
$ws.Range("C5").Value2 = "Ravi0111"
$ws.Range("B5").Value2 = "Adactin.com - New User Registration"

# --- Move the mailto hyperlink from the old F5 position to its new G5 home ---
$ws.Range("F5").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("G5"), "mailto:ravi12@gmail.com")
# Hyperlinks.Add() mints a brand-new cell style; put G5 back on the
# worksheet's existing built-in "Hyperlink" style instead of leaving a
# freshly-added duplicate style behind.
$ws.Range("G5").Style = "Hyperlink"

# --- Column B is now the wide "Expected Title" column (same display width
# as column D, 39 characters -- ColumnWidth is specified in the COM
# "characters" unit, which is slightly smaller than the raw OOXML width) ---
$ws.Columns("B").ColumnWidth = 38.14

# --- Selection moves to the newly edited cell ---
$ws.Range("B5").Select()
